$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "AnyOf" entry as a new row (row 18) below the existing table,
# in column B, matching the plain column style already used for column B.
$ws.Range("B18").Value = "AnyOf"

# Move the selection to C18 to match the author's final cursor position
# (mirrors the updated <selection pane="bottomRight" .../> in the diff).
$ws.Range("C18").Select()
